# Updates cryptos list figures (prices / volume%) and reorders a few
# coin rows, matching the "Updated cryptos list" GitHub Actions commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) and E (Volume(1h)) hold numeric/percentage-looking
# strings that must stay plain text (the sheet stores them as
# inlineStr). Forcing the cell's number format to "@" (Text) before
# assigning the value stops Excel from re-interpreting strings like
# "0.4424" or "10.70" as floating point numbers, then we restore the
# cell's original style so no formatting side effects leak in.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

# Column B (Coin) and C (Link) are plain, unambiguous text - safe to
# assign directly.
function Set-PlainValue($addr, $val) {
    $ws.Range($addr).Value = $val
}


Set-TextValue "D2" "27.754.54"
Set-TextValue "E2" "  +1.37%  "
Set-TextValue "D3" "1.864.94"
Set-TextValue "E3" "  +1.31%  "
Set-TextValue "E4" "  +0.92%  "
Set-TextValue "D5" "323.20"
Set-TextValue "E5" "  +1.31%  "
Set-TextValue "E6" "  +0.83%  "
Set-TextValue "D7" "0.4424"
Set-TextValue "E7" "  +1.67%  "
Set-TextValue "D8" "0.3802"
Set-TextValue "E8" "  +2.23%  "
Set-TextValue "D9" "0.07465"
Set-TextValue "E9" "  +1.71%  "
Set-TextValue "D10" "0.8874"
Set-TextValue "E10" "  +1.32%  "
Set-TextValue "D11" "21.74"
Set-TextValue "E11" "  +1.64%  "
Set-TextValue "D12" "1.873.78"
Set-TextValue "E12" "  -7.44%  "
Set-TextValue "E13" "  +1.34%  "
Set-TextValue "D14" "6.777"
Set-TextValue "E14" "  +1.47%  "
Set-TextValue "D15" "0.07204"
Set-TextValue "E15" "  +0.72%  "
Set-TextValue "D16" "84.36"
Set-TextValue "E16" "  +2.74%  "
Set-TextValue "E17" "  +0.99%  "
Set-TextValue "E18" "  +1.67%  "
Set-TextValue "E19" "  +0.85%  "
Set-TextValue "D20" "15.57"
Set-TextValue "E20" "  +1.19%  "
Set-TextValue "D21" "27.744.93"
Set-TextValue "E21" "  +1.18%  "
Set-TextValue "D22" "5.310"
Set-TextValue "E22" "  +1.12%  "
Set-TextValue "D23" "11.31"
Set-TextValue "E23" "  +1.63%  "
Set-TextValue "D24" "2.090.24"
Set-TextValue "E24" "  -5.06%  "
Set-TextValue "D25" "2.020"
Set-TextValue "E25" "  +5.93%  "
Set-TextValue "D26" "158.21"
Set-TextValue "E26" "  +0.89%  "
Set-TextValue "D27" "18.87"
Set-TextValue "E27" "  +1.80%  "
Set-PlainValue "B28" "InternetComputer(DFINITY)"
Set-PlainValue "C28" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D28" "5.355"
Set-TextValue "E28" "  +1.63%  "
Set-PlainValue "B29" "LidoDAOToken"
Set-PlainValue "C29" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D29" "1.993"
Set-TextValue "E29" "  +3.62%  "
Set-TextValue "D30" "118.96"
Set-TextValue "E30" "  +2.95%  "
Set-TextValue "D31" "0.09059"
Set-TextValue "E31" "  +0.47%  "
Set-TextValue "D32" "1.229"
Set-TextValue "E32" "  +2.45%  "
Set-TextValue "D33" "0.7791"
Set-TextValue "E33" "  +2.64%  "
Set-TextValue "D34" "3.033"
Set-TextValue "E34" "  +6.19%  "
Set-TextValue "D35" "4.594"
Set-TextValue "E35" "  +2.90%  "
Set-TextValue "E36" "  +0.83%  "
Set-TextValue "D37" "1.145"
Set-TextValue "E37" "  -0.40%  "
Set-TextValue "D38" "0.01988"
Set-TextValue "E38" "  +1.64%  "
Set-TextValue "D39" "0.05363"
Set-TextValue "E39" "  +2.16%  "
Set-TextValue "D40" "2.895"
Set-TextValue "E40" "  +3.40%  "
Set-TextValue "D41" "0.5212"
Set-TextValue "E41" "  +1.14%  "
Set-TextValue "D42" "0.1696"
Set-TextValue "E42" "  +2.05%  "
Set-TextValue "D43" "6.912"
Set-TextValue "E43" "  +5.92%  "
Set-TextValue "D44" "8.716"
Set-TextValue "E44" "  +2.90%  "
Set-TextValue "D45" "110.76"
Set-TextValue "E45" "  +2.40%  "
Set-TextValue "D46" "0.06748"
Set-TextValue "E46" "  +7.31%  "
Set-TextValue "D47" "10.70"
Set-TextValue "E47" "  +1.27%  "
Set-PlainValue "B48" "PaxDollar"
Set-PlainValue "C48" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D48" "1.036"
Set-TextValue "E48" "  +0.93%  "
Set-PlainValue "B49" "NEARProtocol"
Set-PlainValue "C49" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D49" "1.716"
Set-TextValue "E49" "  +2.81%  "
Set-PlainValue "B50" "Decentraland"
Set-PlainValue "C50" "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue "D50" "0.4730"
Set-TextValue "E50" "  +2.17%  "
Set-PlainValue "B51" "RenderToken"
Set-PlainValue "C51" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D51" "1.919"
Set-TextValue "E51" "  +1.59%  "
